$d = $word.ActiveDocument

# Locate the paragraph that ends with
# "...the order corresponds to the value of i to be used." and insert a
# brand-new bulleted paragraph right after it (and therefore right before
# the "Angular rotation is degrees per second..." paragraph).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "the order corresponds to the value of i to be used\.") {
        $target = $p
        break
    }
}

if ($null -eq $target) {
    throw "Could not find the anchor paragraph ending in 'the order corresponds to the value of i to be used.'"
}

$null = $target.Range.InsertParagraphAfter()

# The freshly-minted (still empty) paragraph immediately follows $target.
$newPara = $target.Next()

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:numPr>
<w:ilvl w:val="0"/>
<w:numId w:val="1"/>
</w:numPr>
<w:rPr>
<w:sz w:val="28"/>
<w:szCs w:val="28"/>
</w:rPr>
</w:pPr>
<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">If </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>virtualCarSensorStates</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>[i]</w:t></w:r>
<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> is equal to zero, that sensor detects the line.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$null = $newPara.Range.InsertXML($xml)
